# Update the "取得日時" (retrieved-at) timestamp in column A for rows 2-6
# on the "ランサーズ" sheet from 2026-01-29 06:43:37 to 2026-01-29 06:54:34.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-01-29 06:54:34"

for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
